$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after existing data (row 52 is the last used row -> start at 53)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

$newRows = @(
    @("2025-02-18", "sleep", $true, $true),
    @("2025-02-18", "activity", $true, $true),
    @("2025-02-18", "weekly_activity", $false, $false)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
